$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Insert a new bullet "Dialogue will occasionally run over the
# image, and will sometimes trigger when it's not supposed to" immediately
# before the "At the start of the game..." bullet (same ListParagraph /
# numId=2 list). The new bullet also carries the "_GoBack" bookmark, which
# Word relocates here from wherever it last was (the very end of the doc).
# ---------------------------------------------------------------------------

$find1 = $d.Content
$find1.Find.Execute("At the start of the game", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$existingBullet = $find1.Paragraphs.Item(1)
$existingBullet.Range.InsertParagraphBefore()

# Re-resolve the bullet via Find again: objects captured before the
# structural insert keep a stale .Index, so look it up fresh.
$find2 = $d.Content
$find2.Find.Execute("At the start of the game", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$existingBullet2 = $find2.Paragraphs.Item(1)
$newBulletIndex = $existingBullet2.Index - 1
$newBullet = $d.Paragraphs.Item($newBulletIndex)

# Trailing sentinel "X" so the bookmark we add below spans a real
# (non-collapsed) range; a zero-length range right at end-of-paragraph
# confuses Bookmarks.Add in this host, so we add it around the sentinel
# and then delete the sentinel, leaving the bookmark correctly collapsed.
$newBullet.Range.Text = "Dialogue will occasionally run over the image, and will sometimes trigger when it" + [char]8217 + "s not supposed toX"

$newBulletFresh = $d.Paragraphs.Item($newBulletIndex)
$sentinelEnd = $newBulletFresh.Range.End - 1
$sentinelRange = $d.Range($sentinelEnd - 1, $sentinelEnd)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $sentinelRange)

$d.Range($sentinelEnd - 1, $sentinelEnd).Text = ""

# ---------------------------------------------------------------------------
# Change 2: Drop the stale <w:lastRenderedPageBreak/> cached before the
# "Future plan:" heading run (the one before "Known issues and bugs:" is
# untouched). Re-typing the run's text forces the run to be rebuilt
# without the page-break hint while keeping its formatting (sz/szCs).
# ---------------------------------------------------------------------------

$find3 = $d.Content
$find3.Find.Execute("Future plan:", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $find3.Start
$headingEnd = $find3.End

$d.Range($headingStart, $headingEnd).Text = [char]8203
$d.Range($headingStart, $headingStart + 1).Text = "Future plan:"

Write-Output "done"
